$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) updates - force text format to preserve exact formatting
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.839.92'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.734.13'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.26'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.0000'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5146'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2783'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '39.35'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06110'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.746.31'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07023'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.25'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6423'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.530'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '76.82'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.000'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '25.824.27'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.48'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000006621'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.969.70'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.145'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.749'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.122'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '139.64'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.513'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.794'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08321'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.692'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.421'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04482'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.617'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9797'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6124'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.638'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01578'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.951'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9995'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '100.38'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.3822'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7258'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.960'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.05392'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.258'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1122'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '52.91'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '30.06'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.579'

# Volume(1h) column (E) updates
$ws.Range("E2").Value = '  +0.21%  '
$ws.Range("E3").Value = '  -0.63%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("E5").Value = '  -1.75%  '
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("E7").Value = '  +1.26%  '
$ws.Range("E8").Value = '  +4.83%  '
$ws.Range("E9").Value = '  -2.56%  '
$ws.Range("E10").Value = '  -0.58%  '
$ws.Range("E11").Value = '  -0.02%  '
$ws.Range("E12").Value = '  +1.29%  '
$ws.Range("E13").Value = '  -0.29%  '
$ws.Range("E14").Value = '  +3.54%  '
$ws.Range("E15").Value = '  +1.52%  '
$ws.Range("E16").Value = '  -0.75%  '
$ws.Range("E17").Value = '  -0.08%  '
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("E19").Value = '  +0.08%  '
$ws.Range("E20").Value = '  -0.52%  '
$ws.Range("E21").Value = '  +0.60%  '
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("E23").Value = '  +2.39%  '
$ws.Range("E24").Value = '  +6.27%  '
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("E26").Value = '  +2.65%  '
$ws.Range("E28").Value = '  +0.18%  '
$ws.Range("E29").Value = '  +1.47%  '
$ws.Range("E30").Value = '  -0.28%  '
$ws.Range("E31").Value = '  +1.79%  '
$ws.Range("E32").Value = '  +0.88%  '
$ws.Range("E33").Value = '  +1.50%  '
$ws.Range("E34").Value = '  +2.56%  '
$ws.Range("E35").Value = '  -1.15%  '
$ws.Range("E36").Value = '  -1.13%  '
$ws.Range("E37").Value = '  +2.92%  '
$ws.Range("E38").Value = '  +0.70%  '
$ws.Range("E39").Value = '  +1.69%  '
$ws.Range("E40").Value = '  +2.30%  '
$ws.Range("E41").Value = '  -0.10%  '
$ws.Range("E42").Value = '  -1.02%  '
$ws.Range("E43").Value = '  +0.39%  '
$ws.Range("E44").Value = '  -2.57%  '
$ws.Range("E45").Value = '  +1.85%  '
$ws.Range("E47").Value = '  +5.82%  '
$ws.Range("E48").Value = '  +2.89%  '
$ws.Range("E49").Value = '  +1.14%  '
$ws.Range("E50").Value = '  +0.37%  '
$ws.Range("E51").Value = '  +2.74%  '
